$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 16 cell update(s) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3287
$ws.Range("J17").Value = 3287
$ws.Range("L17").Value = 9861
$ws.Range("N17").Value = -10197
$ws.Range("H112").Value = 2099
$ws.Range("J112").Value = 2102.7693
$ws.Range("L112").Value = 6308.3079
$ws.Range("N112").Value = -8524.3079
$ws.Range("H137").Value = 3005.7715
$ws.Range("I137").Value = 2911.889
$ws.Range("K137").Value = 8735.667000000001
$ws.Range("M137").Value = -6185.667000000001
$ws.Range("H138").Value = 5884
$ws.Range("I138").Value = 3248.4
$ws.Range("K138").Value = 9745.200000000001
$ws.Range("M138").Value = -4605.200000000001

# --- Sheet ARM: 49 cell update(s) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 5799.7856
$ws.Range("I28").Value = 3246.077
$ws.Range("K28").Value = 3246.077
$ws.Range("M28").Value = -3054.077
$ws.Range("H32").Value = 11419514
$ws.Range("I32").Value = 12207231
$ws.Range("J32").Value = 6450837.5
$ws.Range("K32").Value = 12207231
$ws.Range("L32").Value = 6450837.5
$ws.Range("M32").Value = -12206944
$ws.Range("N32").Value = -6451411.5
$ws.Range("H45").Value = 3446.5
$ws.Range("J45").Value = 7666.3335
$ws.Range("L45").Value = 7666.3335
$ws.Range("N45").Value = -8420.333500000001
$ws.Range("H61").Value = 4357.1113
$ws.Range("I61").Value = 2544.2942
$ws.Range("K61").Value = 2544.2942
$ws.Range("M61").Value = -2332.2942
$ws.Range("H74").Value = 4416.778
$ws.Range("I74").Value = 4282
$ws.Range("J74").Value = 4686.3335
$ws.Range("K74").Value = 4282
$ws.Range("L74").Value = 4686.3335
$ws.Range("M74").Value = -3408
$ws.Range("N74").Value = -6434.3335
$ws.Range("H77").Value = 4416.778
$ws.Range("I77").Value = 4282
$ws.Range("J77").Value = 4686.3335
$ws.Range("K77").Value = 21410
$ws.Range("L77").Value = 23431.6675
$ws.Range("M77").Value = -17042
$ws.Range("N77").Value = -32167.6675
$ws.Range("H99").Value = 5799.7856
$ws.Range("I99").Value = 3246.077
$ws.Range("K99").Value = 3246.077
$ws.Range("M99").Value = -251.0770000000002
$ws.Range("H122").Value = 3031.0557
$ws.Range("I122").Value = 2305.8215
$ws.Range("K122").Value = 6917.4645
$ws.Range("M122").Value = -4467.4645
$ws.Range("H132").Value = 4125.0537
$ws.Range("I132").Value = 3362.3594
$ws.Range("K132").Value = 10087.0782
$ws.Range("M132").Value = -7557.0782
$ws.Range("H136").Value = 4357.1113
$ws.Range("I136").Value = 2544.2942
$ws.Range("K136").Value = 7632.882599999999
$ws.Range("M136").Value = -5082.882599999999

# --- Sheet BSM: 16 cell update(s) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 33612.25
$ws.Range("J81").Value = 33612.25
$ws.Range("L81").Value = 33612.25
$ws.Range("N81").Value = -35734.25
$ws.Range("H84").Value = 33612.25
$ws.Range("J84").Value = 33612.25
$ws.Range("L84").Value = 100836.75
$ws.Range("N84").Value = -111444.75
$ws.Range("H105").Value = 3199.5264
$ws.Range("I105").Value = 2345.2727
$ws.Range("K105").Value = 2345.2727
$ws.Range("M105").Value = -598.2727
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

# --- Sheet CRP: 32 cell update(s) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4252.978
$ws.Range("I31").Value = 2936.0476
$ws.Range("J31").Value = 5359.2
$ws.Range("K31").Value = 2936.0476
$ws.Range("L31").Value = 5359.2
$ws.Range("M31").Value = -2641.0476
$ws.Range("N31").Value = -5949.2
$ws.Range("H34").Value = 4252.978
$ws.Range("I34").Value = 2936.0476
$ws.Range("J34").Value = 5359.2
$ws.Range("K34").Value = 2936.0476
$ws.Range("L34").Value = 5359.2
$ws.Range("M34").Value = -2734.0476
$ws.Range("N34").Value = -5763.2
$ws.Range("H58").Value = 3132.2778
$ws.Range("I58").Value = 1973.8334
$ws.Range("J58").Value = 5449.1665
$ws.Range("K58").Value = 1973.8334
$ws.Range("L58").Value = 5449.1665
$ws.Range("M58").Value = -1770.8334
$ws.Range("N58").Value = -5855.1665
$ws.Range("H134").Value = 4825.25
$ws.Range("I134").Value = 3109.7
$ws.Range("K134").Value = 9329.099999999999
$ws.Range("M134").Value = -6794.099999999999
$ws.Range("H136").Value = 3132.2778
$ws.Range("I136").Value = 1973.8334
$ws.Range("J136").Value = 5449.1665
$ws.Range("K136").Value = 5921.5002
$ws.Range("L136").Value = 16347.4995
$ws.Range("M136").Value = -3371.5002
$ws.Range("N136").Value = -21447.4995

# --- Sheet CUL: 31 cell update(s) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1955.1428
$ws.Range("J23").Value = 1111.2727
$ws.Range("L23").Value = 3333.8181
$ws.Range("N23").Value = -3803.8181
$ws.Range("H40").Value = 756.6429000000001
$ws.Range("I40").Value = 51
$ws.Range("J40").Value = 2026.8
$ws.Range("K40").Value = 204
$ws.Range("L40").Value = 8107.2
$ws.Range("M40").Value = -135
$ws.Range("N40").Value = -8245.200000000001
$ws.Range("H74").Value = 5416.6665
$ws.Range("H77").Value = 5416.6665
$ws.Range("H86").Value = 457.42856
$ws.Range("J86").Value = 474.5
$ws.Range("L86").Value = 1423.5
$ws.Range("N86").Value = -3795.5
$ws.Range("H89").Value = 457.42856
$ws.Range("J89").Value = 474.5
$ws.Range("L89").Value = 4270.5
$ws.Range("N89").Value = -16126.5
$ws.Range("I122").Value = 1999
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 17991
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -15541
$ws.Range("N122").Value = -49900
$ws.Range("H131").Value = 5436.25
$ws.Range("J131").Value = 6048.8
$ws.Range("L131").Value = 18146.4
$ws.Range("N131").Value = -28226.4

# --- Sheet GSM: 4 cell update(s) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 61666.5
$ws.Range("J135").Value = 61666.5
$ws.Range("L135").Value = 61666.5
$ws.Range("N135").Value = -71806.5

# --- Sheet LTW: 21 cell update(s) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2421.4
$ws.Range("I22").Value = 1550.5
$ws.Range("J22").Value = 3002
$ws.Range("K22").Value = 1550.5
$ws.Range("L22").Value = 3002
$ws.Range("M22").Value = -1255.5
$ws.Range("N22").Value = -3592
$ws.Range("H27").Value = 2421.4
$ws.Range("I27").Value = 1550.5
$ws.Range("J27").Value = 3002
$ws.Range("K27").Value = 1550.5
$ws.Range("L27").Value = 3002
$ws.Range("M27").Value = -1443.5
$ws.Range("N27").Value = -3216
$ws.Range("H122").Value = 7946.923
$ws.Range("I122").Value = 5126
$ws.Range("J122").Value = 9200.666999999999
$ws.Range("K122").Value = 15378
$ws.Range("L122").Value = 27602.001
$ws.Range("M122").Value = -12928
$ws.Range("N122").Value = -32502.001

# --- Sheet WVR: 11 cell update(s) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1996
$ws.Range("I107").Value = 2003.1111
$ws.Range("K107").Value = 6009.3333
$ws.Range("M107").Value = -4089.3333
$ws.Range("H136").Value = 3533.2246
$ws.Range("I136").Value = 3038.9355
$ws.Range("J136").Value = 4384.5
$ws.Range("K136").Value = 9116.806500000001
$ws.Range("L136").Value = 13153.5
$ws.Range("M136").Value = -6566.806500000001
$ws.Range("N136").Value = -18253.5
